$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.7
$ws.Range("D3").Value = -7.7
$ws.Range("D5").Value = -7.961
$ws.Range("E7").Value = 13.078
$ws.Range("A9").Value = -20.775
$ws.Range("E9").Value = 12.734
$ws.Range("D11").Value = -8.15
$ws.Range("D12").Value = -8.088999999999999
$ws.Range("A13").Value = -21.99
$ws.Range("A16").Value = -20.825
$ws.Range("A18").Value = -21.834
$ws.Range("A20").Value = -21.763
$ws.Range("D21").Value = -7.813999999999998
$ws.Range("E21").Value = 13.272
